$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").Value = -6.150299999999995
$ws.Range("D4").Value = -7.872100000000001
$ws.Range("E6").Value = 12.42230000000001
$ws.Range("D7").Value = -8.261399999999997
$ws.Range("E7").Value = 12.3513
$ws.Range("D8").Value = -8.483799999999995
$ws.Range("E8").Value = 12.5858
$ws.Range("A11").Value = -21.94390000000003
$ws.Range("A12").Value = -22.77640000000002
$ws.Range("D12").Value = -8.46010000000001
$ws.Range("D14").Value = -8.672800000000002
$ws.Range("A15").Value = -21.44610000000002
$ws.Range("E19").Value = 12.8894
$ws.Range("E21").Value = 12.60819999999999
$ws.Range("D22").Value = -8.007899999999999
$ws.Range("E24").Value = 12.96769999999999
$ws.Range("E25").Value = 13.17250000000001
